$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158, shifting existing rows 158:170 down to 159:171
$ws.Rows.Item(158).Insert(-4121)

# Populate the newly inserted row 158 with the new weekly record
$ws.Cells.Item(158, 1).Value = 5
$ws.Cells.Item(158, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(158, 3).Value = "Maule"
$ws.Cells.Item(158, 4).Value = 44461
$ws.Cells.Item(158, 5).Value = 7
$ws.Cells.Item(158, 6).Value = 100112009
$ws.Cells.Item(158, 7).Value = "Acelga"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 500
$ws.Cells.Item(158, 11).Value = 2300
$ws.Cells.Item(158, 12).Value = 2300
$ws.Cells.Item(158, 13).Value = 2300
$ws.Cells.Item(158, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(158, 15).Value = "Región del Maule"
$ws.Cells.Item(158, 16).Value = 575
$ws.Cells.Item(158, 17).Value = 4
$ws.Cells.Item(158, 18).Value = "Hortaliza"

"done"
